$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 8.799974000000001
$ws.Range("N2").Value = 26.399922
$ws.Range("O2").Value = 0.1502220177021807
$ws.Range("P2").Value = 0.1502220177021807
$ws.Range("Q2").Value = 3.715868221266
$ws.Range("R2").Value = 33.442813991394
$ws.Range("S2").Value = 0.1502220177021807
$ws.Range("T2").Value = 0.1502220177021807

# Row 3
$ws.Range("O3").Value = 0.4499951903206205
$ws.Range("P3").Value = 0.4499951903206205
$ws.Range("S3").Value = 0.4499951903206205
$ws.Range("T3").Value = 0.4499951903206205

# Row 4
$ws.Range("M4").Value = 23.07309566666667
$ws.Range("N4").Value = 69.21928700000001
$ws.Range("O4").Value = 0.3938746848208995
$ws.Range("P4").Value = 0.3938746848208995
$ws.Range("Q4").Value = 9.742822303111002
$ws.Range("R4").Value = 87.68540072799901
$ws.Range("S4").Value = 0.3938746848208995
$ws.Range("T4").Value = 0.3938746848208995

# Row 5
$ws.Range("M5").Value = 0.3460956666666666
$ws.Range("N5").Value = 1.038287
$ws.Range("O5").Value = 0.005908107156299329
$ws.Range("P5").Value = 0.00590810715629933
$ws.Range("Q5").Value = 0.146142010111
$ws.Range("R5").Value = 1.315278090999
$ws.Range("S5").Value = 0.005908107156299329
$ws.Range("T5").Value = 0.00590810715629933
